$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (runs, balls, fours, sixes) -> 0, 2, 0, 0
$ws.Range("C2:F2").NumberFormat = "@"
$ws.Cells.Item(2, 3).Value = "0"
$ws.Cells.Item(2, 4).Value = "2"
$ws.Cells.Item(2, 5).Value = "0"
$ws.Cells.Item(2, 6).Value = "0"

# Row 4 (runs, balls, fours, sixes) -> 34, 36, 1, 2
$ws.Range("C4:F4").NumberFormat = "@"
$ws.Cells.Item(4, 3).Value = "34"
$ws.Cells.Item(4, 4).Value = "36"
$ws.Cells.Item(4, 5).Value = "1"
$ws.Cells.Item(4, 6).Value = "2"
